$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Priority column (E) for rows 4-7 changes from "low" to "ht" in both locale sheets
$wsZhCn.Range("E4:E7").Value2 = "ht"
$wsDeDe.Range("E4:E7").Value2 = "ht"

# Latest Handoff Datetime (H) for zh-cn rows 4-7: bump generated timestamp
$wsZhCn.Range("H4:H7").Value2 = "2016-09-04 12:35:01"

# Latest HO Xliff Generate Date shared text used by Overview!G4:G7 and de-de!H4:H7
$wsOverview.Range("G4:G7").Value2 = "2016-09-04 12:35:11"
$wsDeDe.Range("H4:H7").Value2 = "2016-09-04 12:35:11"
